# TPB test.pptx - "Add files via upload" edit
#
# Slide 1 ("TPB Test"): the dimensions textbox (TextBox 6) gets a new
# trailing line "u = 5mm/s" appended after "Notch width = 4 mm". The
# shape uses <a:spAutoFit/> so PowerPoint grows its height automatically
# to fit the extra line - we just change the text and let autofit do it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$dimsShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text -like "*Notch width*") {
            $dimsShape = $candidate
        }
    }
}

if ($dimsShape -ne $null) {
    $tr = $dimsShape.TextFrame.TextRange
    # Add a new paragraph at the end of the existing text.
    $tr.Text = $tr.Text + "`r" + "u = 5mm/s"
}
